$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newNote = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.7 = 36304.35 pesos`n✅ 36304.35 pesos = 8.66 = 953.1 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newNote

# --- Update the rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 115
$ws2.Range("O10").Value = 4175
$ws2.Range("N12").Value = 4190
$ws2.Range("O12").Value = 110
